$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(548).Insert()

$ws.Range("A548").Value = 5
$ws.Range("B548").Value = "Macroferia Regional de Talca"
$ws.Range("C548").Value = "Maule"
$ws.Range("D548").Value = 45124
$ws.Range("E548").Value = 7
$ws.Range("F548").Value = 100114014
$ws.Range("G548").Value = "Betarraga"
$ws.Range("H548").Value = "Sin especificar"
$ws.Range("I548").Value = "Primera"
$ws.Range("J548").Value = 5000
$ws.Range("K548").Value = 600
$ws.Range("L548").Value = 600
$ws.Range("M548").Value = 600
$ws.Range("N548").Value = '$/paquete 5 unidades'
$ws.Range("O548").Value = 'Región del Maule'
$ws.Range("P548").Value = 120
$ws.Range("Q548").Value = 5
$ws.Range("R548").Value = "Hortaliza"

Write-Output "done"
